$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 75.18000000000001
$ws.Range("C2").Value = 'The candidate has a score of 75.18 which indicates a high aptness for the SDE Intern position. The candidate''s experience in working on projects involving MERN stack development, ReactJS, JavaScript, and Python aligns well with the required skills for the job. Additionally, their experience in implementing machine learning models and data analysis demonstrates their expertise in relevant technologies.'
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 'The applicant demonstrates a strong willingness to come to Japan by expressing interest in the country''s leading position in AI, the desire to be a part of Japan''s work culture, and the interest in learning a new language. Additionally, the applicant''s career plan aligns with the company''s focus on AI/ML. The applicant''s willingness to work in Japan is further reinforced by their openness to working in different countries. In terms of personality, the applicant displays strengths in teamwork, problem-solving, and adaptability. The experience shared reflects the ability to overcome challenges through collaboration and adaptability, aligning with the company''s desired soft skills. The applicant also has a unique background in both engineering and graphic design, which sets them apart and adds value to the company''s projects. The applicant''s demonstrated willingness, career plan, and alignment with the company''s needs contribute to the overall positive score.'

# Row 3
$ws.Range("B3").Value = 82.98
$ws.Range("C3").Value = 'The applicant has a high score because they have experience working with the required technologies (MongoDB, ReactJS, JavaScript, Web Development, NodeJS) in their projects. They have successfully developed web/mobile applications, utilized WebSockets and WebRTC APIs, and displayed data in easy-to-understand graphs. Their diverse project experience demonstrates their ability to handle different aspects of the job description.'
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 'The applicant demonstrates a strong willingness to come to Japan, expressing interest in the country''s AI industry and work culture. They also show adaptability to other cultures and a desire to learn a new language. Their career plan aligns with the company''s focus on AI/ML. In terms of personality, the applicant exhibits strengths in teamwork, adaptability, and problem-solving, which are crucial soft skills for the company. They have also reflected on their weaknesses and taken steps to address them, showing a willingness to improve and grow. Overall, the applicant presents a strong fit for the company and shows a genuine interest in working in Japan.'

# Row 4
$ws.Range("B4").Value = 75.28
$ws.Range("C4").Value = 'The candidate has a high score of 75.28 as their projects show significant experience in web and mobile application development using MongoDB, ReactJS, JavaScript, and NodeJS. Additionally, their projects demonstrate proficiency in computer vision, image processing, and machine learning using Python, PyTorch, Tensorflow, and Sklearn. These skills closely align with the required skills for the SDE Intern position, making the candidate well-suited for the job.'
$ws.Range("E4").Value = 'The applicant demonstrates a strong willingness to come to Japan, citing the country''s leading position in AI and expressing a genuine interest in the work culture. The career plan and the willingness to work in different countries, including Japan, show flexibility and adaptability. The applicant''s strengths in teamwork and problem-solving, along with the ability to bridge technical and creative aspects, make them a good fit for the company. However, the preference not to work alone and the need to improve skills could be perceived as minor concerns. Overall, the applicant shows potential and a positive attitude, aligning well with the company''s values and requirements.'
